# Generate Report for Handback
# - Updates the "Ready for handoff" status for 3ef7e9ba-... to "Handback transform failed"
# - Adds detailed "Error Detail" messages for the zh-cn and de-de handback sheets
# - Widens the "Error Detail" column to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text (row for 3ef7e9ba-...) everywhere it appears: the Overview
# sheet (zh-cn/de-de status columns) and the "Status" column on each per-locale sheet.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Set the "Error Detail" column (P) values on the per-locale sheets for the 3ef7e9ba-... row
$wsZhCn.Range("P3").Value = "Handback file name: yu2zwmdy.olh is different with handoff file name: 3ef7e9ba-13d3-4123-b80f-121cfd22aa52.c5abf32b9b2d4c603aeb7f1d1124a769da518ee9.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: yu2zwmdy.olh is different with handoff file name: 3ef7e9ba-13d3-4123-b80f-121cfd22aa52.c5abf32b9b2d4c603aeb7f1d1124a769da518ee9.de-de."

# Widen the "Error Detail" column (P) on both locale sheets to fit the new content.
# (ColumnWidth is specified in characters; Excel stores column width in OOXML after
# adding ~5/6 of a character of internal cell padding, so subtract that back out here
# so the saved width attribute comes out to exactly 40.)
$wsZhCn.Range("P1:P1").ColumnWidth = 40 - 5/6
$wsDeDe.Range("P1:P1").ColumnWidth = 40 - 5/6
